$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, reordered/updated triage data (Component, Bugs, Features, Total)
$data = @(
    @("Menu", 16, 5, 21),
    @("Combobox", 7, 5, 12),
    @("Nav", 9, 1, 10),
    @("DataGrid", 4, 5, 9),
    @("Popover", 6, 2, 8),
    @("Dropdown", 4, 3, 7),
    @("Tree", 7, 0, 7),
    @("Dialog", 4, 2, 6),
    @("Toolbar", 3, 1, 4),
    @("Table", 4, 0, 4),
    @("Drawer", 4, 0, 4),
    @("Tooltip", 3, 1, 4),
    @("Virtualizer", 4, 0, 4),
    @("Skeleton", 2, 1, 3),
    @("TeachingPopover", 3, 0, 3),
    @("Calendar Compat", 3, 0, 3),
    @("TagPicker", 2, 1, 3),
    @("List", 2, 0, 2),
    @("Portal", 2, 0, 2),
    @("MessageBar", 2, 0, 2),
    @("Slider", 0, 2, 2),
    @("DatePicker", 0, 2, 2),
    @("Tabs", 2, 0, 2),
    @("FluentProvider", 0, 2, 2),
    @("Switch", 2, 0, 2),
    @("Button", 0, 1, 1),
    @("InfoLabel", 1, 0, 1),
    @("Checkbox", 1, 0, 1),
    @("AvatarGroup", 1, 0, 1),
    @("DatePickerCompat", 1, 0, 1),
    @("Toast", 1, 0, 1),
    @("Input", 0, 1, 1),
    @("Accordion", 1, 0, 1),
    @("Spinner", 0, 1, 1),
    @("Badge", 1, 0, 1),
    @("Popup", 1, 0, 1),
    @("Avatar", 1, 0, 1),
    @("Card", 1, 0, 1),
    @("Image", 0, 1, 1),
    @("FocusTrapZone", 1, 0, 1),
    @("Label", 1, 0, 1),
    @("SplitButton", 1, 0, 1),
    @("MenuItem", 1, 0, 1),
    @("Tag", 0, 0, 0),
    @("Pickers", 0, 0, 0),
    @("Carousel", 0, 0, 0),
    @("Rating", 0, 0, 0),
    @("SpinButton", 0, 0, 0),
    @("Coachmark", 0, 0, 0),
    @("Segment", 0, 0, 0),
    @("Calendar", 0, 0, 0),
    @("Keytip", 0, 0, 0),
    @("ColorPicker", 0, 0, 0)
)

# Clear out the previous data region below the header (old range went to row 50)
$ws.Range("A2:D50").Clear()

$rowIndex = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]
    $ws.Cells.Item($rowIndex, 3).Value = $entry[2]
    $ws.Cells.Item($rowIndex, 4).Value = $entry[3]
    $rowIndex++
}
